$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete now-unused trailing rows (43 and 42) so remaining rows shift up / dimension shrinks
$ws.Rows.Item(43).Delete()
$ws.Rows.Item(42).Delete()

# Update changed cell values to match target dataset
$ws.Range("B2").Value = 670.5865089999995
$ws.Range("I2").Value = 2.930399
$ws.Range("O2").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "C2"
$ws.Range("I3").Value = 2.940398999999996
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = "C2"
$ws.Range("O3").Value = 0
$ws.Range("E4").Value = 1
$ws.Range("I4").Value = 1.9701
$ws.Range("K4").Value = 1
$ws.Range("O4").Value = 1.940399
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = "C2"
$ws.Range("I5").Value = 1.95999898989899
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = "C2"
$ws.Range("O5").Value = 2.58719866666666
$ws.Range("E6").Value = 2
$ws.Range("I6").Value = 1.291099501
$ws.Range("K6").Value = 2
$ws.Range("O6").Value = 0
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = "C2"
$ws.Range("I7").Value = 3.57719866666667
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = "C2"
$ws.Range("O7").Value = 2.587198666666669
$ws.Range("E8").Value = 3
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 3
$ws.Range("O8").Value = 1.940399
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = "C2"
$ws.Range("G9").Value = "P1"
$ws.Range("H9").Value = "CPU1"
$ws.Range("I9").Value = 2.940399000000001
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = "C2"
$ws.Range("O9").Value = 2.587198666666663
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = "C1"
$ws.Range("G10").Value = "P1"
$ws.Range("H10").Value = "CPU1"
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = 4
$ws.Range("O10").Value = 0
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = "C2"
$ws.Range("G11").Value = "P1"
$ws.Range("H11").Value = "CPU1"
$ws.Range("I11").Value = 2.940399000000001
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 4
$ws.Range("L11").Value = "C2"
$ws.Range("O11").Value = 2.587198666666669
$ws.Range("K12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = "C2"
$ws.Range("O13").Value = 2.561326679999993
$ws.Range("K14").Value = 1
$ws.Range("O14").Value = 0
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = "C2"
$ws.Range("O15").Value = 0
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = 2
$ws.Range("O16").Value = 1.950399
$ws.Range("J17").Value = 2
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = "C2"
$ws.Range("O17").Value = 2.56132668
$ws.Range("J18").Value = 2
$ws.Range("K18").Value = 3
$ws.Range("O18").Value = 0
$ws.Range("J19").Value = 2
$ws.Range("L19").Value = "C2"
$ws.Range("O19").Value = 2.561326680000002
$ws.Range("J20").Value = 2
$ws.Range("O20").Value = 0
$ws.Range("J21").Value = 2
$ws.Range("K21").Value = 4
$ws.Range("L21").Value = "C2"
$ws.Range("O21").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("J23").Value = 3
$ws.Range("L23").Value = "C2"
$ws.Range("O23").Value = 0
$ws.Range("J24").Value = 3
$ws.Range("O24").Value = 0
$ws.Range("J25").Value = 3
$ws.Range("K25").Value = 1
$ws.Range("L25").Value = "C2"
$ws.Range("O25").Value = 0
$ws.Range("J26").Value = 3
$ws.Range("K26").Value = 2
$ws.Range("O26").Value = 0
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = 2
$ws.Range("L27").Value = "C2"
$ws.Range("O27").Value = 0
$ws.Range("J28").Value = 3
$ws.Range("K28").Value = 3
$ws.Range("O28").Value = 0
$ws.Range("J29").Value = 3
$ws.Range("K29").Value = 3
$ws.Range("L29").Value = "C2"
$ws.Range("O29").Value = 0
$ws.Range("J30").Value = 3
$ws.Range("K30").Value = 4
$ws.Range("O30").Value = 0
$ws.Range("J31").Value = 3
$ws.Range("K31").Value = 4
$ws.Range("L31").Value = "C2"
$ws.Range("O31").Value = 0
$ws.Range("J32").Value = 4
$ws.Range("K32").Value = 0
$ws.Range("O32").Value = 0
$ws.Range("J33").Value = 4
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = "C2"
$ws.Range("O33").Value = 0
$ws.Range("J34").Value = 4
$ws.Range("K34").Value = 1
$ws.Range("O34").Value = 0
$ws.Range("J35").Value = 4
$ws.Range("K35").Value = 1
$ws.Range("L35").Value = "C2"
$ws.Range("O35").Value = 0
$ws.Range("J36").Value = 4
$ws.Range("K36").Value = 2
$ws.Range("O36").Value = 0
$ws.Range("J37").Value = 4
$ws.Range("K37").Value = 2
$ws.Range("L37").Value = "C2"
$ws.Range("O37").Value = 0
$ws.Range("J38").Value = 4
$ws.Range("K38").Value = 3
$ws.Range("O38").Value = 0
$ws.Range("J39").Value = 4
$ws.Range("K39").Value = 3
$ws.Range("L39").Value = "C2"
$ws.Range("O39").Value = 0
$ws.Range("J40").Value = 4
$ws.Range("K40").Value = 4
$ws.Range("O40").Value = 0
$ws.Range("J41").Value = 4
$ws.Range("L41").Value = "C2"
$ws.Range("O41").Value = 0
